## Fix the "error-duplicates" example-weights sheet:
## Row 5 was an accidental duplicate of row 2 (same "from" value, different
## weight). Shift the from-column values up one slot and correct the
## duplicate row's weight, then fold the previously-excluded row into the
## Table1 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table1 only covered A1:D4 before; grow it to include row 5 (also updates
# the table's autoFilter range to match).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D5"))

# De-duplicate: shift the "from" values for rows 3-5 up one position and
# give row 5 the corrected (non-duplicate) weight.
$ws.Range("B3").Value = "Driver capabilities and limitations"
$ws.Range("B4").Value = "Road Transportation System capabilities and limitations"
$ws.Range("B5").Value = "Vehicle capacity"
$ws.Range("D5").Value = 0.75

# Column C no longer needs to be as wide once recalculated.
$ws.Columns.Item(3).ColumnWidth = 16.83

# Leave the selection on the first data row instead of the stray D6 cell.
$ws.Range("A2").Select()
